# Automatische test-sync: 2025-08-03 02:12:50
# Adds a new log row (row 3) to the "Logs" sheet, widens the conditional
# formatting ranges that highlight that row's columns, and bumps the
# "Intern verzoek / Actie voor medewerker" counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- New row of data (row 3) -------------------------------------------------
$ws.Range("A3").Value = "Kun jij dit even regelen?"
$ws.Range("B3").Value = "mailmind.test@zohomail.eu"
$ws.Range("C3").Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Range("D3").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E3").Value = "Beste klant, `nBedankt voor uw e-mail. Om u zo goed mogelijk van dienst te kunnen zijn, ontvang ik graag meer informatie over wat u precies geregeld wilt hebben. Kunt u mij wat meer details geven, zodat ik u beter kan helpen? Alvast bedankt voor uw medewerking.`nMet vriendelijke groet, `n[Naam]`nE-mailassistent"
$ws.Range("F3").Value = "2025-08-03 02:12:12"
$ws.Range("G3").Value = "Ja"
$ws.Range("H3").Value = "Nee"
$ws.Range("I3").Value = "Ja"
$ws.Range("J3").Value = "Nee"

# --- Extend conditional formatting sqref of each column to include row 3 ----
foreach ($col in @("D", "G", "H", "I", "J")) {
    $newRange = $ws.Range("$($col)2:$($col)3")
    $fcs = $ws.Range("$($col)2").FormatConditions
    $count = $fcs.Count()
    for ($i = 1; $i -le $count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard counter bump ---------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 2
